$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 111.46154
$ws.Range("I12").Value = 111.46154
$ws.Range("K12").Value = 111.46154
$ws.Range("M12").Value = 58.53846
$ws.Range("H17").Value = 968662.75
$ws.Range("J17").Value = 989692.1
$ws.Range("L17").Value = 2969076.3
$ws.Range("N17").Value = -2969412.3
$ws.Range("H33").Value = 994.4
$ws.Range("I33").Value = 634.8570999999999
$ws.Range("K33").Value = 634.8570999999999
$ws.Range("M33").Value = -405.8570999999999
$ws.Range("H40").Value = 2821.3635
$ws.Range("J40").Value = 3348
$ws.Range("L40").Value = 3348
$ws.Range("N40").Value = -3698
$ws.Range("H88").Value = 1030.1428
$ws.Range("I88").Value = 760.75
$ws.Range("K88").Value = 760.75
$ws.Range("M88").Value = -354.75
$ws.Range("H91").Value = 1030.1428
$ws.Range("I91").Value = 760.75
$ws.Range("K91").Value = 760.75
$ws.Range("M91").Value = 643.25
$ws.Range("H92").Value = 654.58826
$ws.Range("I92").Value = 626.8214
$ws.Range("J92").Value = 784.1667
$ws.Range("K92").Value = 626.8214
$ws.Range("L92").Value = 784.1667
$ws.Range("M92").Value = 621.1786
$ws.Range("N92").Value = -3280.1667
$ws.Range("H95").Value = 48000
$ws.Range("I95").Value = 34000
$ws.Range("J95").Value = 55000
$ws.Range("K95").Value = 34000
$ws.Range("L95").Value = 55000
$ws.Range("M95").Value = -31254
$ws.Range("N95").Value = -60492
$ws.Range("H106").Value = 4944496.5
$ws.Range("I106").Value = 5373148.5
$ws.Range("J106").Value = 14997.5
$ws.Range("K106").Value = 5373148.5
$ws.Range("L106").Value = 14997.5
$ws.Range("M106").Value = -5372517.5
$ws.Range("N106").Value = -16259.5
$ws.Range("H107").Value = 4650.727
$ws.Range("I107").Value = 4475.9443
$ws.Range("J107").Value = 5437.25
$ws.Range("K107").Value = 4475.9443
$ws.Range("L107").Value = 5437.25
$ws.Range("M107").Value = -2555.9443
$ws.Range("N107").Value = -9277.25
$ws.Range("H115").Value = 1397.5
$ws.Range("I115").Value = 621.25
$ws.Range("K115").Value = 1863.75
$ws.Range("M115").Value = -296.75
$ws.Range("H137").Value = 8312.596
$ws.Range("I137").Value = 13162.305
$ws.Range("J137").Value = 2441.8948
$ws.Range("K137").Value = 39486.915
$ws.Range("L137").Value = 7325.6844
$ws.Range("M137").Value = -36936.915
$ws.Range("N137").Value = -12425.6844
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6743.7974
$ws.Range("I32").Value = 6465.366
$ws.Range("K32").Value = 6465.366
$ws.Range("M32").Value = -6178.366
$ws.Range("H74").Value = 7770.1763
$ws.Range("I74").Value = 11888
$ws.Range("J74").Value = 3137.625
$ws.Range("K74").Value = 11888
$ws.Range("L74").Value = 3137.625
$ws.Range("M74").Value = -11014
$ws.Range("N74").Value = -4885.625
$ws.Range("H77").Value = 7770.1763
$ws.Range("I77").Value = 11888
$ws.Range("J77").Value = 3137.625
$ws.Range("K77").Value = 59440
$ws.Range("L77").Value = 15688.125
$ws.Range("M77").Value = -55072
$ws.Range("N77").Value = -24424.125
$ws.Range("H92").Value = 366683330
$ws.Range("J92").Value = 366683330
$ws.Range("L92").Value = 366683330
$ws.Range("N92").Value = -366688322
$ws.Range("H132").Value = 2849.682
$ws.Range("I132").Value = 1930.1936
$ws.Range("J132").Value = 5042.3076
$ws.Range("K132").Value = 5790.5808
$ws.Range("L132").Value = 15126.9228
$ws.Range("M132").Value = -3260.5808
$ws.Range("N132").Value = -20186.9228
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H43").Value = 208855
$ws.Range("J43").Value = 208855
$ws.Range("L43").Value = 208855
$ws.Range("N43").Value = -209217
$ws.Range("H86").Value = 8405.143
$ws.Range("I86").Value = 9393
$ws.Range("J86").Value = 4783
$ws.Range("K86").Value = 9393
$ws.Range("L86").Value = 4783
$ws.Range("M86").Value = -8270
$ws.Range("N86").Value = -7029
$ws.Range("H89").Value = 8405.143
$ws.Range("I89").Value = 9393
$ws.Range("J89").Value = 4783
$ws.Range("K89").Value = 46965
$ws.Range("L89").Value = 23915
$ws.Range("M89").Value = -41349
$ws.Range("N89").Value = -35147
$ws.Range("H105").Value = 127777.78
$ws.Range("I105").Value = 222500
$ws.Range("J105").Value = 9375
$ws.Range("K105").Value = 222500
$ws.Range("L105").Value = 9375
$ws.Range("M105").Value = -220753
$ws.Range("N105").Value = -12869
$ws.Range("H107").Value = 2643.6155
$ws.Range("I107").Value = 2697.25
$ws.Range("K107").Value = 2697.25
$ws.Range("M107").Value = -777.25
$ws.Range("H134").Value = 6880.0605
$ws.Range("I134").Value = 6905.1035
$ws.Range("K134").Value = 20715.3105
$ws.Range("M134").Value = -18180.3105
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6229.9346
$ws.Range("I31").Value = 6586.9355
$ws.Range("K31").Value = 6586.9355
$ws.Range("M31").Value = -6291.9355
$ws.Range("H34").Value = 6229.9346
$ws.Range("I34").Value = 6586.9355
$ws.Range("K34").Value = 6586.9355
$ws.Range("M34").Value = -6384.9355
$ws.Range("H122").Value = 10464.462
$ws.Range("I122").Value = 12988.9
$ws.Range("K122").Value = 38966.7
$ws.Range("M122").Value = -36516.7
$ws.Range("H134").Value = 4728.7095
$ws.Range("I134").Value = 5275.4443
$ws.Range("K134").Value = 15826.3329
$ws.Range("M134").Value = -13291.3329
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H137").Value = 4020.7727
$ws.Range("J137").Value = 15983.333
$ws.Range("L137").Value = 47949.999
$ws.Range("N137").Value = -58149.999
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 10584.1
$ws.Range("I70").Value = 12016.5
$ws.Range("K70").Value = 12016.5
$ws.Range("M70").Value = -11746.5
$ws.Range("H73").Value = 10584.1
$ws.Range("I73").Value = 12016.5
$ws.Range("K73").Value = 12016.5
$ws.Range("M73").Value = -11080.5
$ws.Range("H97").Value = 6327.346
$ws.Range("I97").Value = 6229.625
$ws.Range("K97").Value = 6229.625
$ws.Range("M97").Value = -5733.625
$ws.Range("H102").Value = 6651.375
$ws.Range("I102").Value = 8035.222
$ws.Range("J102").Value = 2499.8333
$ws.Range("K102").Value = 8035.222
$ws.Range("L102").Value = 2499.8333
$ws.Range("M102").Value = -6413.222
$ws.Range("N102").Value = -5743.8333
$ws.Range("H132").Value = 4350.913
$ws.Range("J132").Value = 2866
$ws.Range("L132").Value = 8598
$ws.Range("N132").Value = -13658
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 14707.667
$ws.Range("I7").Value = 17685.822
$ws.Range("K7").Value = 17685.822
$ws.Range("M7").Value = -17573.822
$ws.Range("H40").Value = 21968.1
$ws.Range("I40").Value = 25003.316
$ws.Range("J40").Value = 16725.455
$ws.Range("K40").Value = 25003.316
$ws.Range("L40").Value = 16725.455
$ws.Range("M40").Value = -24867.316
$ws.Range("N40").Value = -16997.455
$ws.Range("H93").Value = 3848.2942
$ws.Range("I93").Value = 5055.1665
$ws.Range("J93").Value = 951.8
$ws.Range("K93").Value = 5055.1665
$ws.Range("L93").Value = 951.8
$ws.Range("M93").Value = -3807.1665
$ws.Range("N93").Value = -3447.8
$ws.Range("H126").Value = 14707.667
$ws.Range("I126").Value = 17685.822
$ws.Range("K126").Value = 53057.466
$ws.Range("M126").Value = -50587.466
$ws.Range("H136").Value = 5741.5835
$ws.Range("I136").Value = 3955.2222
$ws.Range("J136").Value = 6813.4
$ws.Range("K136").Value = 11865.6666
$ws.Range("L136").Value = 20440.2
$ws.Range("M136").Value = -9315.6666
$ws.Range("N136").Value = -25540.2
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 482469.62
$ws.Range("J62").Value = 14938.25
$ws.Range("L62").Value = 14938.25
$ws.Range("N62").Value = -16186.25
$ws.Range("H65").Value = 482469.62
$ws.Range("J65").Value = 14938.25
$ws.Range("L65").Value = 74691.25
$ws.Range("N65").Value = -80931.25
$ws.Range("H122").Value = 17844.828
$ws.Range("I122").Value = 1898.45
$ws.Range("K122").Value = 5695.35
$ws.Range("M122").Value = -3245.35
$ws.Range("H132").Value = 8106.27
$ws.Range("I132").Value = 9618.762000000001
$ws.Range("K132").Value = 28856.286
$ws.Range("M132").Value = -26326.286
$ws.Range("H136").Value = 488972.16
$ws.Range("I136").Value = 675533.3
$ws.Range("J136").Value = 12204.777
$ws.Range("K136").Value = 2026599.9
$ws.Range("L136").Value = 36614.331
$ws.Range("M136").Value = -2024049.9
$ws.Range("N136").Value = -41714.331
